$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 107
$ws.Range("H107").Value = 1102.4375
$ws.Range("I107").Value = 991.0833
$ws.Range("J107").Value = 1436.5
$ws.Range("K107").Value = 991.0833
$ws.Range("L107").Value = 1436.5
$ws.Range("M107").Value = 928.9167
$ws.Range("N107").Value = -5276.5
# Row 132
$ws.Range("H132").Value = 3157.0312
$ws.Range("I132").Value = 3439.3928
$ws.Range("J132").Value = 1180.5
$ws.Range("K132").Value = 10318.1784
$ws.Range("L132").Value = 3541.5
$ws.Range("M132").Value = -7788.178400000001
$ws.Range("N132").Value = -8601.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1294.6364
$ws.Range("I2").Value = 1294.6364
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1294.6364
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1181.6364
$ws.Range("N2").ClearContents()
# Row 110
$ws.Range("H110").Value = 993.6667
$ws.Range("I110").Value = 841.1429000000001
$ws.Range("K110").Value = 841.1429000000001
$ws.Range("M110").Value = 1203.8571
# Row 114
$ws.Range("H114").Value = 34665.332
$ws.Range("J114").Value = 34665.332
$ws.Range("L114").Value = 34665.332
$ws.Range("N114").Value = -43343.332
# Row 116
$ws.Range("H116").Value = 1294.6364
$ws.Range("I116").Value = 1294.6364
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1294.6364
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 999.3635999999999
$ws.Range("N116").ClearContents()
# Row 132
$ws.Range("H132").Value = 10925.574
$ws.Range("I132").Value = 1539
$ws.Range("J132").Value = 31355.176
$ws.Range("K132").Value = 4617
$ws.Range("L132").Value = 94065.52799999999
$ws.Range("M132").Value = -2087
$ws.Range("N132").Value = -99125.52799999999
# Row 133
$ws.Range("H133").Value = 62498
$ws.Range("J133").Value = 62498
$ws.Range("L133").Value = 62498
$ws.Range("N133").Value = -67558

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1294.6364
$ws.Range("I3").Value = 1294.6364
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1294.6364
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1180.6364
$ws.Range("N3").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 271.91666
$ws.Range("I7").Value = 16.5
$ws.Range("J7").Value = 527.3333
$ws.Range("K7").Value = 16.5
$ws.Range("L7").Value = 527.3333
$ws.Range("M7").Value = 96.5
$ws.Range("N7").Value = -753.3333

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1258.55
$ws.Range("I5").Value = 973.38464
$ws.Range("J5").Value = 1788.1428
$ws.Range("K5").Value = 2920.15392
$ws.Range("L5").Value = 5364.428400000001
$ws.Range("M5").Value = -2808.15392
$ws.Range("N5").Value = -5588.428400000001
# Row 14
$ws.Range("H14").Value = 989.8
$ws.Range("I14").Value = 989.8
$ws.Range("K14").Value = 2969.4
$ws.Range("M14").Value = -2796.4
# Row 25
$ws.Range("H25").Value = 1649.75
$ws.Range("I25").Value = 1649.75
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 4949.25
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -4780.25
$ws.Range("N25").ClearContents()
# Row 30
$ws.Range("H30").Value = 1649.75
$ws.Range("I30").Value = 1649.75
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 4949.25
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -4847.25
$ws.Range("N30").ClearContents()
# Row 131
$ws.Range("H131").Value = 674.91
$ws.Range("J131").Value = 699.5714
$ws.Range("L131").Value = 2098.7142
$ws.Range("N131").Value = -12178.7142
# Row 135
$ws.Range("H135").Value = 1258.55
$ws.Range("I135").Value = 973.38464
$ws.Range("J135").Value = 1788.1428
$ws.Range("K135").Value = 8760.46176
$ws.Range("L135").Value = 16093.2852
$ws.Range("M135").Value = -6225.46176
$ws.Range("N135").Value = -21163.2852

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 9153.154
$ws.Range("I113").Value = 9499.25
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 9499.25
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -7329.25
$ws.Range("N113").Value = -9340
# Row 126
$ws.Range("H126").Value = 3524.4
$ws.Range("I126").Value = 3581.4707
$ws.Range("K126").Value = 10744.4121
$ws.Range("M126").Value = -8274.4121
# Row 141
$ws.Range("H141").Value = 56811
$ws.Range("J141").Value = 56811
$ws.Range("L141").Value = 56811
$ws.Range("N141").Value = -67171

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2631.8572
$ws.Range("I22").Value = 3174.75
$ws.Range("J22").Value = 894.6
$ws.Range("K22").Value = 3174.75
$ws.Range("L22").Value = 894.6
$ws.Range("M22").Value = -2879.75
$ws.Range("N22").Value = -1484.6
# Row 27
$ws.Range("H27").Value = 2631.8572
$ws.Range("I27").Value = 3174.75
$ws.Range("J27").Value = 894.6
$ws.Range("K27").Value = 3174.75
$ws.Range("L27").Value = 894.6
$ws.Range("M27").Value = -3067.75
$ws.Range("N27").Value = -1108.6
# Row 61
$ws.Range("H61").Value = 5183
$ws.Range("I61").Value = 2103.75
$ws.Range("J61").Value = 17500
$ws.Range("K61").Value = 2103.75
$ws.Range("L61").Value = 17500
$ws.Range("M61").Value = -1901.75
$ws.Range("N61").Value = -17904
# Row 113
$ws.Range("H113").Value = 5183
$ws.Range("I113").Value = 2103.75
$ws.Range("J113").Value = 17500
$ws.Range("K113").Value = 2103.75
$ws.Range("L113").Value = 17500
$ws.Range("M113").Value = 66.25
$ws.Range("N113").Value = -21840
# Row 136
$ws.Range("H136").Value = 2160.0952
$ws.Range("I136").Value = 1634
$ws.Range("K136").Value = 4902
$ws.Range("M136").Value = -2352

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 724.75
$ws.Range("I122").Value = 724.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2174.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 275.75
$ws.Range("N122").ClearContents()
# Row 136
$ws.Range("H136").Value = 25809406
$ws.Range("I136").Value = 32259272
$ws.Range("K136").Value = 96777816
$ws.Range("M136").Value = -96775266
# Row 140
$ws.Range("H140").Value = 45425
$ws.Range("J140").Value = 45425
$ws.Range("L140").Value = 45425
$ws.Range("N140").Value = -55785
# Row 141
$ws.Range("H141").Value = 73250
$ws.Range("J141").Value = 73250
$ws.Range("L141").Value = 73250
$ws.Range("N141").Value = -83610
